$wb = $excel.ActiveWorkbook

# Add a new worksheet "Test4" after the last existing sheet (this makes it
# the active sheet, which automatically updates workbookView's activeTab
# and moves tabSelected from the previously active sheet to this one).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Test4"

# Header row
$ws.Range("D8").Value = "A"
$ws.Range("E8").Value = "B"
$ws.Range("F8").Value = "C"
$ws.Range("G8").Value = "D"

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "a"
$ws.Range("F9").Value = "z"

$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "b"
$ws.Range("F10").Value = "y"
$ws.Range("G10").Value = 1

$ws.Range("E11").Value = "c"
$ws.Range("F11").Value = "x"
$ws.Range("G11").Value = 2

$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "d"
$ws.Range("F12").Value = "w"
$ws.Range("G12").Value = 3

$ws.Range("D13").Value = 4
$ws.Range("E13").Value = "e"
$ws.Range("G13").Value = 4

$ws.Range("D14").Value = 5
$ws.Range("E14").Value = "f"
$ws.Range("F14").Value = "v"
$ws.Range("G14").Value = 5

$ws.Range("D15").Value = 6
$ws.Range("E15").Value = "g"
$ws.Range("F15").Value = "u"

$ws.Range("E16").Value = "h"

# Set the selection on the new sheet to C5, matching the target workbook
$ws.Range("C5").Select() | Out-Null
